$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "matt_karting"
$ws.Range("B2").Value = "z"
$ws.Range("C2").Value = "z"
$ws.Range("D2").Value = "z"
$ws.Range("E2").Value = "z"

$ws.Range("A3").Value = "palmer_jrr"
$ws.Range("B3").Value = "vd"
$ws.Range("C3").Value = "scv"
$ws.Range("D3").Value = "scsv"
$ws.Range("E3").Value = "vd"
